$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: plain text "*" (not number-like, so no quote-prefix style needed)
$ws.Range("A3").Value = "*"
$ws.Range("B3").Value = "*"
$ws.Range("C3").Value = "*"
$ws.Range("D3").Value = "*"
$ws.Range("E3").Value = "*"
$ws.Range("H3").Value = "*"

# Row 3 has an explicit row height of 14 (customHeight)
$ws.Rows("3:3").RowHeight = 14

# Row 4: B4 is a genuine number (123) styled with a quote-prefix (text-like) format.
# Build the quote-prefix style via a scratch cell, then copy its format onto B4.
$ws.Range("Z1").Value = "'999"
$ws.Range("B4").Value = 123
$ws.Range("Z1").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Range("C4").Value = "*"

# Row 5
$ws.Range("D5").Value = "*"
$ws.Range("G5").Value = "*"

# Row 6
$ws.Range("E6").Value = "*"

# Row 7
$ws.Range("C7").Value = "*"

# Final selection
$ws.Range("B10").Select()
